# T2 - Tabla puntuación: mark several checklist cells with "X",
# set a one-off underline style on D17, and move the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1 (rows 9-12): mark all four items as done ---
$ws.Range("D9").Value  = "X"
$ws.Range("D10").Value = "X"
$ws.Range("D11").Value = "X"
$ws.Range("D12").Value = "X"

# --- Block 2 (rows 17-19): mark D/E cells as done ---
$ws.Range("D17").Value = "X"
$ws.Range("E17").Value = "X"
$ws.Range("D18").Value = "X"
$ws.Range("E18").Value = "X"
$ws.Range("D19").Value = "X"
$ws.Range("E19").Value = "X"

# D17 gets a distinct underline style (new font/cellXf) vs the rest
$ws.Range("D17").Font.Underline = $true

# --- Block 3 (rows 22-28): mark a couple of items as done ---
$ws.Range("D25").Value = "X"
$ws.Range("E25").Value = "X"
$ws.Range("D27").Value = "X"

# --- Row-height normalisation: several rows shrink back to the
#     sheet default (14.25) once re-laid-out, the two wrapped-text
#     header rows (14 & 21) settle on a slightly shorter height ---
$ws.Rows("2:4").EntireRow.AutoFit()
$ws.Rows("6").EntireRow.AutoFit()
$ws.Rows("10").EntireRow.AutoFit()
$ws.Rows("12").EntireRow.AutoFit()
$ws.Rows("17:19").EntireRow.AutoFit()
$ws.Rows("22:26").EntireRow.AutoFit()
$ws.Rows("28").EntireRow.AutoFit()

$ws.Rows("14").RowHeight = 25.5
$ws.Rows("21").RowHeight = 25.5

# --- Scroll / selection: the view now sits lower on the sheet ---
$ws.Range("E33").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 2
